# Updated reconstruction tools with revision data
#
# Updates the "It meets current standards" row (row 12) for tools AU and RA,
# and the "Provides synonyms for metabolites and reactions" row (row 23) for
# tools AU and ME, with revised descriptions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "Inputs in SBML L3 and outputs in SBML L3 with some features missing (MIRIAM compliant CV annotations and SBML Groups). Be aware that slightly different output networks can be obtained when using input networks in SBML L2 format (in SBML L3, there will be metabolites explicitily stored in the compartment boundary). 1 out of 5 points are discounted."

$ws.Range("J12").Value = "outputs in SBML L3 with FBC annotations, SBML Groups and MIRIAM compliant CV annotations. Be aware that networks created with RAVEN have to be exported to SBML using the specific functions of RAVEN (not COBRA functions as a regular COBRA user would expect) because otherwise there will be no MIRIAM annotations in the SBML files. 0 out of 5 points is discounted"

$ws.Range("D23").Value = "It provides synonyms for reactions but they can only be found in files in PADMET format and not in the SBML files. Synonyms for metabolites are missing. 2 out of 5 points are discounted"

$ws.Range("G23").Value = "It only provides identifiers from the reference database. No synonyms for other databases are provided. 5 out of 5 points are discounted"

# Row heights adjust automatically for wrapped text in Excel, but set them
# explicitly to match the reviewed layout.
$ws.Rows.Item(7).RowHeight = 54.9
$ws.Rows.Item(12).RowHeight = 115.2
$ws.Rows.Item(23).RowHeight = 57.6

# Restore the view state (scroll position / active cell) recorded on save.
$ws.Range("C20").Select()
